$d = $word.ActiveDocument

$d.Content.Find.Execute("As everything, an AI must be tested. There might happen that a model makes an error, like human does. Sometimes we are unable to recognize a person or we see a dog instead of a cat. If a PC ", $true, $false, $false, $false, $false, $true, 1, $false, "Today ", 2)
$d.Content.Find.Execute("isn’t", $true, $false, $false, $false, $false, $true, 1, $false, "I’ve", 2)
$d.Content.Find.Execute(" always right about the classification of a dog, that might be fatal for an autopilot by not seeing a pedestrian. If human can classify an emotion in the face correctly at the 70%, an AI cannot go a lot further", $true, $false, $false, $false, $false, $true, 1, $false, " tried an Image Recognition model. The first exercise was to train and test a model to recognize dogs, cats and dolphins. The second one asked to train and test a model to recognize two people with different clothes. The last one was the same as the first one but the background changes. I noticed that the AI couldn’t work really well if the 2 people changed their dresses as well as if the background changed with animals", 2)
$d.Content.Find.Execute("If we have a lot of data, we should always split it. The 4/5 in training data and 1/5 in test set. This is useful in order that the AI can work well on new data.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

Write-Output $d.Content.Text
